# Applies the two changes captured by the commit:
#   1. Slide 6's table switches to the built-in table style
#      {7F733C68-96EC-4217-BC29-8EAA429FD874}.
#   2. The presentation's theme colour scheme (the "Integral" palette that
#      theme2.xml/SlideMaster.Theme carries) is swapped for the stock
#      "Office Theme" palette that used to live in theme1.xml.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{7F733C68-96EC-4217-BC29-8EAA429FD874}", $false)
    }
}

# --- 2. Theme colour scheme: Integral -> Office Theme ---------------------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

function Set-ThemeColor {
    param($Index, $R, $G, $B)
    $packed = $R + ($G * 256) + ($B * 65536)
    $colorScheme.Item($Index).RGB = $packed
}

Set-ThemeColor 1  0x00 0x00 0x00   # dk1
Set-ThemeColor 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor 12 0x95 0x4F 0x72   # folHlink
